$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-NewsLink {
    param(
        [string]$RangeAddress,
        [string]$Url
    )
    $rng = $ws.Range($RangeAddress)
    # Create the hyperlink first (anchored to the top-left cell of the range);
    # then (re)write the cell value(s) and apply the Hyperlink style so every
    # cell in a merged/grouped range shares the same "visited link" look
    # without the engine cloning a near-duplicate style record.
    $ws.Hyperlinks.Add($rng, $Url)
    $rng.Value = $Url
    $rng.Style = "Hyperlink"
}

# D36 - Cedar Meats original article (15 March)
Add-NewsLink "D36" "https://www2.health.vic.gov.au/about/media-centre/MediaReleases/more-covid19-cases-confirmed-victoria-15-march"

# D7 - 20 May update
Add-NewsLink "D7" "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-20-may-2020"

# D8 - ABC article, Cedar Meats / Kyabram health workers
Add-NewsLink "D8" "https://www.abc.net.au/news/2020-05-20/coronavirus-victoria-cedar-meats-outbreak-kyabram-health-workers/12266102"

# D9 & D10 - 19 May update
Add-NewsLink "D9" "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-19-may-2020"
Add-NewsLink "D10" "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-19-may-2020"

# D11 - 18 May update
Add-NewsLink "D11" "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-18-may-2020"

# D12:D23 - 18 May update, shared across the whole block of rows
Add-NewsLink "D12:D23" "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-18-may-2020"

# D24 - 16 May update
Add-NewsLink "D24" "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-16-may"

# D26 - 17 May update
Add-NewsLink "D26" "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-17-may-2020"

# D25 - ABC article, 15 May McDonald's / Cedar Meats
Add-NewsLink "D25" "https://www.abc.net.au/news/2020-05-15/victoria-coronavirus-cases-rise-by-20-mcdonalds-cedar-meats/12251762"
